# Remove the last six slides from the deck (slides 18-23), matching the
# author's edit that dropped the "A Note About Transformations of
# Variables...", "College Student Heights", "More properties of Linear
# Transformations", "Sampling Distribution" (x2) and "Margin of Error"
# slides from the end of the presentation.

$p = $ppt.ActivePresentation

# Delete from the highest index down to 18 so indices of slides we still
# want to keep never shift out from under us.
for ($i = $p.Slides.Count; $i -ge 18; $i--) {
    $p.Slides.Item($i).Delete()
}
